# Update the dSF column (F) values for several rows to reflect the
# repulled data / recalculated means, per commit "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = -3
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 4
$ws.Range("F10").Value = -4
$ws.Range("F12").Value = -7
$ws.Range("F14").Value = -8
$ws.Range("F16").Value = 3
